$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Novembro")

$ws.Range("A3").Value = "25/10/2023"
$ws.Range("B3").NumberFormat = '#,##0.00 €; [Red]-#,##0.00 €'
$ws.Range("E3").Value = "Income"

$ws.Range("A5").Value = "26/10/2023"
$ws.Range("B5").NumberFormat = '#,##0.00 €; [Red]-#,##0.00 €'
$ws.Range("E5").Value = "Income"

$ws.Range("A18").Value = "14/11/2023"
$ws.Range("B18").NumberFormat = '#,##0.00 €; [Red]-#,##0.00 €'
$ws.Range("E18").Value = "Income"
